# Iteration5: Adding built code
# Adds 5 new defect-log rows (27-31) to Sheet1, matching rows logged on
# 2019-05-01 (serial date 43586), plus the trailing page-setup / selection
# tweaks that came along with the save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New data rows -----------------------------------------------------
# Columns: A=Date, B=No., C=Inject, D=Type, E=Remove, F=Fix Time,
#          G=(=B<row>), H=Description
$rows = @(
    @{ Row = 27; No = 21; Inject = 20; FixTime = 2; Height = 30;
       Desc = '[Vue warn]: Error in data(): "TypeError: _model.GameOne is not a constructor" - Need to export class' },
    @{ Row = 28; No = 22; Inject = 60; FixTime = 1; Height = 45;
       Desc = '[Vue warn]: Invalid prop: type check failed for prop "guess". Expected Number with value NaN, got String with value "Cheating?". - String check (typeof/instanceof String)' },
    @{ Row = 29; No = 23; Inject = 80; FixTime = 1; Height = 30;
       Desc = "[Vue warn]: Error in v-on handler: ""TypeError: Cannot read property 'toString' of undefined"" - Extra else statement needed to return string" },
    @{ Row = 30; No = 24; Inject = 70; FixTime = 1; Height = 30;
       Desc = "[Vue warn]: Error in render: ""TypeError: Cannot read property 'toString' of undefined"" - Needed super(newGuess, gameNumber)" },
    @{ Row = 31; No = 25; Inject = 60; FixTime = 5; Height = 45;
       Desc = 'Endless Loop. Crashed browser. When check if number had been used before, if all numbers had been used then it would crash the browser. - Added check to prevent this' }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Copy the date-formatted cell above so the new A cell picks up the
    # same number format (style 18) instead of minting a duplicate style.
    $ws.Range("A26").Copy()
    $ws.Range("A$rowNum").PasteSpecial(-4122)
    $ws.Range("A$rowNum").Value = 43586

    $ws.Range("B$rowNum").Value = $r.No
    $ws.Range("C$rowNum").Value = $r.Inject
    $ws.Range("D$rowNum").Value = "Code"
    $ws.Range("E$rowNum").Value = "Compile"
    $ws.Range("F$rowNum").Value = $r.FixTime
    $ws.Range("G$rowNum").Formula = "=B$rowNum"
    $ws.Range("H$rowNum").Value = $r.Desc

    $ws.Rows.Item($rowNum).RowHeight = $r.Height
}

# --- Sheet view / selection ---------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("F31").Select() | Out-Null

# --- Page setup -----------------------------------------------------------
$ws.PageSetup.Orientation = 1
